$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "Translations"

# Row 1 (header row) - shift headers right by one and add new "Entity Id" header
$ws.Cells.Item(1,1).Value = "Entity Id"
$ws.Cells.Item(1,2).Value = "Type"
$ws.Cells.Item(1,3).Value = "Index"
$ws.Cells.Item(1,4).Value = "Original"
$ws.Cells.Item(1,5).Value = "Translation"

# Row 2
$ws.Cells.Item(2,1).Value = "AAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAA"
$ws.Cells.Item(2,2).Value = "Title"
$ws.Cells.Item(2,3).ClearContents()
$ws.Cells.Item(2,4).Value = "Orig"

# Row 3
$ws.Cells.Item(3,1).Value = "AAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAA"
$ws.Cells.Item(3,2).Value = "ValidationMessage"
$ws.Cells.Item(3,3).Value = 1
$ws.Cells.Item(3,4).Value = "Orig"
$ws.Cells.Item(3,5).Value = "validation message"

# Row 4
$ws.Cells.Item(4,1).Value = "AAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAA"
$ws.Cells.Item(4,2).Value = "Instruction"
$ws.Cells.Item(4,3).ClearContents()
$ws.Cells.Item(4,4).Value = "Orig"

# Row 5
$ws.Cells.Item(5,1).Value = "AAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAA"
$ws.Cells.Item(5,2).Value = "OptionTitle"
$ws.Cells.Item(5,3).Value = 2
$ws.Cells.Item(5,4).Value = "Orig"
$ws.Cells.Item(5,5).Value = "option"

# Resize columns to fit the new content (widths tuned to match target layout)
$ws.Columns.Item(1).ColumnWidth = 42.5
$ws.Columns.Item(2).ColumnWidth = 17.333333333333332
$ws.Columns.Item(3).ColumnWidth = 5.166666666666667
$ws.Columns.Item(5).ColumnWidth = 39.0

# Update selection to reflect new active cell
$ws.Range("E6").Select() | Out-Null
